# Update "想去人数" (interested-people count) figures by +1 on both the
# "展览" sheet and the "全部类型" sheet, matching the regenerated data
# snapshot published to gh-pages.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1449
$ws1.Range("F6").Value  = 1727
$ws1.Range("F14").Value = 76
$ws1.Range("F20").Value = 4525
$ws1.Range("F27").Value = 2024

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1449
$ws4.Range("F6").Value  = 1727
$ws4.Range("F14").Value = 76
$ws4.Range("F20").Value = 4525
$ws4.Range("F29").Value = 2024
